$d = $word.ActiveDocument

# Grab the full flat-OPC package representation (includes every part: document.xml,
# styles.xml, numbering.xml, ...) so we can surgically edit the numbering + styles
# parts and feed the whole thing back through InsertXML in one shot.
$full = $d.Content.WordOpenXML

$oldNumBlock = @'
<w:num w:numId="56" w16cid:durableId="424033570"><w:abstractNumId w:val="2"/><w:lvlOverride w:ilvl="0"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="1"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="2"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="3"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="4"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="5"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="6"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="7"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="8"><w:startOverride w:val="1"/></w:lvlOverride></w:num>
'@

$newNumBlock = @'
<w:num w:numId="56" w16cid:durableId="424033570"><w:abstractNumId w:val="2"/><w:lvlOverride w:ilvl="0"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="1"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="2"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="3"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="4"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="5"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="6"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="7"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="8"><w:startOverride w:val="1"/></w:lvlOverride></w:num><w:num w:numId="57" w16cid:durableId="2056199994"><w:abstractNumId w:val="1"/></w:num><w:num w:numId="58" w16cid:durableId="77946263"><w:abstractNumId w:val="3"/><w:lvlOverride w:ilvl="0"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="1"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="2"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="3"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="4"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="5"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="6"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="7"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="8"><w:startOverride w:val="1"/></w:lvlOverride></w:num><w:num w:numId="59" w16cid:durableId="2102750283"><w:abstractNumId w:val="2"/><w:lvlOverride w:ilvl="0"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="1"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="2"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="3"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="4"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="5"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="6"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="7"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="8"><w:startOverride w:val="1"/></w:lvlOverride></w:num><w:num w:numId="60" w16cid:durableId="69893764"><w:abstractNumId w:val="2"/><w:lvlOverride w:ilvl="0"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="1"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="2"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="3"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="4"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="5"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="6"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="7"><w:startOverride w:val="1"/></w:lvlOverride><w:lvlOverride w:ilvl="8"><w:startOverride w:val="1"/></w:lvlOverride></w:num>
'@

$oldStyleSnip = @'
<w:style w:type="paragraph" w:customStyle="1" w:styleId="FirstParagraph"><w:name w:val="First Paragraph"/><w:basedOn w:val="Textoindependiente"/><w:next w:val="Textoindependiente"/><w:qFormat/></w:style>
'@

$newStyleSnip = @'
<w:style w:type="paragraph" w:customStyle="1" w:styleId="FirstParagraph"><w:name w:val="First Paragraph"/><w:basedOn w:val="Textoindependiente"/><w:next w:val="Textoindependiente"/><w:qFormat/><w:rsid w:val="0022278C"/><w:pPr><w:ind w:firstLine="0"/></w:pPr></w:style>
'@

if ($full.IndexOf($oldNumBlock) -lt 0) {
    throw "numbering.xml anchor block not found"
}
if ($full.IndexOf($oldStyleSnip) -lt 0) {
    throw "styles.xml anchor block not found"
}

$full = $full.Replace($oldNumBlock, $newNumBlock)
$full = $full.Replace($oldStyleSnip, $newStyleSnip)

$d.Content.InsertXML($full)

Write-Host "Applied numbering + style updates"
